$d = $word.ActiveDocument

# Replace the inline picture (Range-based Parking Provision Standards diagram)
# with a hyperlink whose display text is the image's URL.
$shp = $d.InlineShapes.Item(1)
$pos = $shp.Range.Start
$shp.Delete()

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Range_Based_Car_Parking_Standard.jpg?h=100%25&w=100%25"
$r = $d.Range($pos, $pos)
$d.Hyperlinks.Add($r, $url, $null, $null, $url) | Out-Null
